$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture "LM35 ADC" (currently F5) into C6, and "MOTOR" (currently E6) into F2:F4 ---
# Copy F5 ("LM35 ADC", style s6) -> C6 (value + format)
$ws.Range("F5").Copy()
$ws.Range("C6").PasteSpecial(-4163)
$ws.Range("F5").Copy()
$ws.Range("C6").PasteSpecial(-4122)

# Copy E6 ("MOTOR", style s8) -> F2, F3, F4, F5 (value + format).
# Do this before E6/E7 get cleared below.
foreach ($dest in "F2", "F3", "F4", "F5") {
    $ws.Range("E6").Copy()
    $ws.Range($dest).PasteSpecial(-4163)
    $ws.Range("E6").Copy()
    $ws.Range($dest).PasteSpecial(-4122)
}

# --- Step 2: D2:D5 become plain black cells (style s5, same as B2/B3/B4/B5) ---
$ws.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("G9").PasteSpecial(-4122)

# --- Step 3: clear E6/E7 to bold-black empty cells (new style matching fontId2/fillId5/border/center) ---
$ws.Range("E6").Value = ""
$ws.Range("E6").ClearFormats()
$ws.Range("B2").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Font.Bold = $true

$ws.Range("E7").Value = ""
$ws.Range("E7").ClearFormats()
$ws.Range("B2").Copy()
$ws.Range("E7").PasteSpecial(-4122)
$ws.Range("E7").Font.Bold = $true

# --- Step 4: clear E8/E9 to plain bordered empty cells (new style: border only, no fill/align) ---
$ws.Range("E8").Value = ""
$ws.Range("E8").ClearFormats()
$ws.Range("E8").Borders.Weight = 2
$ws.Range("E8").Borders.LineStyle = 1

$ws.Range("E9").Value = ""
$ws.Range("E9").ClearFormats()
$ws.Range("E9").Borders.Weight = 2
$ws.Range("E9").Borders.LineStyle = 1

# --- Step 5: update selection to match target ---
$ws.Range("E23").Select()
